# Scheduled market-data refresh: push updated Universalis price snapshots
# into the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) of
# each job sheet. Generated by the runner from the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 479.8
$ws.Range("I31").Value = 479.8
$ws.Range("K31").Value = 1439.4
$ws.Range("M31").Value = -1209.4
$ws.Range("H39").Value = 194.5
$ws.Range("I39").Value = 126
$ws.Range("J39").Value = 400
$ws.Range("K39").Value = 378
$ws.Range("L39").Value = 1200
$ws.Range("M39").Value = -82
$ws.Range("N39").Value = -1792
$ws.Range("H98").Value = 1213.125
$ws.Range("I98").Value = 1172.1428
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 1172.1428
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 325.8571999999999
$ws.Range("N98").Value = -4496
$ws.Range("H111").Value = 83334880
$ws.Range("I111").Value = 142858780
$ws.Range("J111").Value = 1399.2
$ws.Range("K111").Value = 428576340
$ws.Range("L111").Value = 4197.6
$ws.Range("M111").Value = -428573273
$ws.Range("N111").Value = -10331.6
$ws.Range("H122").Value = 1213.125
$ws.Range("I122").Value = 1172.1428
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3516.4284
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1066.4284
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 50000104
$ws.Range("I3").Value = 205
$ws.Range("K3").Value = 205
$ws.Range("M3").Value = -90
$ws.Range("H45").Value = 2127.1904
$ws.Range("I45").Value = 1642.6364
$ws.Range("K45").Value = 1642.6364
$ws.Range("M45").Value = -1265.6364
$ws.Range("H61").Value = 3314.0908
$ws.Range("I61").Value = 2744.182
$ws.Range("J61").Value = 3884
$ws.Range("K61").Value = 2744.182
$ws.Range("L61").Value = 3884
$ws.Range("M61").Value = -2532.182
$ws.Range("N61").Value = -4308
$ws.Range("H88").Value = 1978.9286
$ws.Range("I88").Value = 2008.75
$ws.Range("J88").Value = 1800
$ws.Range("K88").Value = 2008.75
$ws.Range("L88").Value = 1800
$ws.Range("M88").Value = -1602.75
$ws.Range("N88").Value = -2612
$ws.Range("H91").Value = 1978.9286
$ws.Range("I91").Value = 2008.75
$ws.Range("J91").Value = 1800
$ws.Range("K91").Value = 2008.75
$ws.Range("L91").Value = 1800
$ws.Range("M91").Value = -604.75
$ws.Range("N91").Value = -4608
$ws.Range("H122").Value = 1280095.1
$ws.Range("I122").Value = 1471834.5
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 4415503.5
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -4413053.5
$ws.Range("N122").Value = -10400.0002
$ws.Range("H132").Value = 3010.258
$ws.Range("I132").Value = 2598.476
$ws.Range("J132").Value = 3875
$ws.Range("K132").Value = 7795.428
$ws.Range("L132").Value = 11625
$ws.Range("M132").Value = -5265.428
$ws.Range("N132").Value = -16685
$ws.Range("H136").Value = 3314.0908
$ws.Range("I136").Value = 2744.182
$ws.Range("J136").Value = 3884
$ws.Range("K136").Value = 8232.545999999998
$ws.Range("L136").Value = 11652
$ws.Range("M136").Value = -5682.545999999998
$ws.Range("N136").Value = -16752

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1950
$ws.Range("I5").Value = 1950
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1950
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1837
$ws.Range("N5").ClearContents()
$ws.Range("H86").Value = 2581.2693
$ws.Range("I86").Value = 2385.875
$ws.Range("J86").Value = 2893.9
$ws.Range("K86").Value = 2385.875
$ws.Range("L86").Value = 2893.9
$ws.Range("M86").Value = -1262.875
$ws.Range("N86").Value = -5139.9
$ws.Range("H89").Value = 2581.2693
$ws.Range("I89").Value = 2385.875
$ws.Range("J89").Value = 2893.9
$ws.Range("K89").Value = 11929.375
$ws.Range("L89").Value = 14469.5
$ws.Range("M89").Value = -6313.375
$ws.Range("N89").Value = -25701.5
$ws.Range("H134").Value = 2127.743
$ws.Range("I134").Value = 2060.7083
$ws.Range("J134").Value = 2274
$ws.Range("K134").Value = 6182.124899999999
$ws.Range("L134").Value = 6822
$ws.Range("M134").Value = -3647.124899999999
$ws.Range("N134").Value = -11892

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2135.5757
$ws.Range("J31").Value = 6098.7144
$ws.Range("L31").Value = 6098.7144
$ws.Range("N31").Value = -6688.7144
$ws.Range("H34").Value = 2135.5757
$ws.Range("J34").Value = 6098.7144
$ws.Range("L34").Value = 6098.7144
$ws.Range("N34").Value = -6502.7144
$ws.Range("H62").Value = 2002358
$ws.Range("I62").Value = 2502197.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2502197.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2501573.5
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2002358
$ws.Range("I65").Value = 2502197.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 12510987.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -12507867.5
$ws.Range("N65").Value = -21240
$ws.Range("H94").Value = 4525.5713
$ws.Range("J94").Value = 6844.7144
$ws.Range("L94").Value = 6844.7144
$ws.Range("N94").Value = -7746.7144
$ws.Range("H122").Value = 833.02856
$ws.Range("I122").Value = 788.6799999999999
$ws.Range("J122").Value = 943.9
$ws.Range("K122").Value = 2366.04
$ws.Range("L122").Value = 2831.7
$ws.Range("M122").Value = 83.96000000000004
$ws.Range("N122").Value = -7731.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30052.107
$ws.Range("I4").Value = 433.33334
$ws.Range("J4").Value = 38129.953
$ws.Range("K4").Value = 1300.00002
$ws.Range("L4").Value = 114389.859
$ws.Range("M4").Value = -1188.00002
$ws.Range("N4").Value = -114613.859
$ws.Range("H5").Value = 1030.1111
$ws.Range("J5").Value = 668.3333
$ws.Range("L5").Value = 2004.9999
$ws.Range("N5").Value = -2228.9999
$ws.Range("H122").Value = 711.0909
$ws.Range("J122").Value = 799.5
$ws.Range("L122").Value = 7195.5
$ws.Range("N122").Value = -12095.5
$ws.Range("H129").Value = 14921.5
$ws.Range("I129").Value = 3367.8
$ws.Range("J129").Value = 34177.668
$ws.Range("K129").Value = 10103.4
$ws.Range("L129").Value = 102533.004
$ws.Range("M129").Value = -5103.400000000001
$ws.Range("N129").Value = -112533.004
$ws.Range("H131").Value = 850.6
$ws.Range("J131").Value = 935.85364
$ws.Range("L131").Value = 2807.56092
$ws.Range("N131").Value = -12887.56092
$ws.Range("H135").Value = 1030.1111
$ws.Range("J135").Value = 668.3333
$ws.Range("L135").Value = 6014.9997
$ws.Range("N135").Value = -11084.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11116133
$ws.Range("I122").Value = 25004100
$ws.Range("J122").Value = 5760
$ws.Range("K122").Value = 75012300
$ws.Range("L122").Value = 17280
$ws.Range("M122").Value = -75009850
$ws.Range("N122").Value = -22180
$ws.Range("H126").Value = 2838.875
$ws.Range("I126").Value = 2878.6155
$ws.Range("J126").Value = 2666.6667
$ws.Range("K126").Value = 8635.8465
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("M126").Value = -6165.8465
$ws.Range("N126").Value = -12940.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 34682.25
$ws.Range("J128").Value = 34682.25
$ws.Range("L128").Value = 34682.25
$ws.Range("N128").Value = -44642.25
$ws.Range("H132").Value = 19582.6
$ws.Range("I132").Value = 27173.9
$ws.Range("K132").Value = 81521.70000000001
$ws.Range("M132").Value = -78991.70000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 10000000
$ws.Range("I8").Value = 10000000
$ws.Range("K8").Value = 10000000
$ws.Range("M8").Value = -9999860
$ws.Range("H54").Value = 15380
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 18966.666
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 18966.666
$ws.Range("M54").Value = -9480
$ws.Range("N54").Value = -20006.666
$ws.Range("H81").Value = 1400
$ws.Range("I81").Value = 1200
$ws.Range("K81").Value = 2400
$ws.Range("M81").Value = -1339
$ws.Range("H84").Value = 1400
$ws.Range("I84").Value = 1200
$ws.Range("K84").Value = 12000
$ws.Range("M84").Value = -6696
$ws.Range("H126").Value = 975.375
$ws.Range("I126").Value = 767.25
$ws.Range("J126").Value = 1599.75
$ws.Range("K126").Value = 2301.75
$ws.Range("L126").Value = 4799.25
$ws.Range("M126").Value = 168.25
$ws.Range("N126").Value = -9739.25
